$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '57.930.56' }
    @{ Cell = 'E2'; Value = '  -0.25%  ' }
    @{ Cell = 'D3'; Value = '2.428.76' }
    @{ Cell = 'E3'; Value = '  -0.97%  ' }
    @{ Cell = 'E4'; Value = '  +0.03%  ' }
    @{ Cell = 'D5'; Value = '509.76' }
    @{ Cell = 'E5'; Value = '  -2.87%  ' }
    @{ Cell = 'D6'; Value = '133.06' }
    @{ Cell = 'E6'; Value = '  +1.64%  ' }
    @{ Cell = 'D7'; Value = '0.997' }
    @{ Cell = 'E7'; Value = '  -0.29%  ' }
    @{ Cell = 'E8'; Value = '  -1.57%  ' }
    @{ Cell = 'D9'; Value = '2.448.92' }
    @{ Cell = 'E9'; Value = '  -0.39%  ' }
    @{ Cell = 'D10'; Value = '0.0976' }
    @{ Cell = 'E10'; Value = '  -0.17%  ' }
    @{ Cell = 'E11'; Value = '  -0.97%  ' }
    @{ Cell = 'D13'; Value = '4.60' }
    @{ Cell = 'E13'; Value = '  -7.66%  ' }
    @{ Cell = 'D14'; Value = '2.843.93' }
    @{ Cell = 'E14'; Value = '  -1.57%  ' }
    @{ Cell = 'D15'; Value = '57.424.80' }
    @{ Cell = 'E15'; Value = '  -1.02%  ' }
    @{ Cell = 'D16'; Value = '21.80' }
    @{ Cell = 'E16'; Value = '  +0.48%  ' }
    @{ Cell = 'D17'; Value = '0.0000133' }
    @{ Cell = 'E17'; Value = '  +0.40%  ' }
    @{ Cell = 'D18'; Value = '2.430.39' }
    @{ Cell = 'E18'; Value = '  -1.02%  ' }
    @{ Cell = 'D19'; Value = '10.28' }
    @{ Cell = 'E19'; Value = '  -1.66%  ' }
    @{ Cell = 'D20'; Value = '4.11' }
    @{ Cell = 'E20'; Value = '  +0.02%  ' }
    @{ Cell = 'D21'; Value = '314.48' }
    @{ Cell = 'E21'; Value = '  -0.21%  ' }
    @{ Cell = 'D22'; Value = '6.40' }
    @{ Cell = 'E22'; Value = '  +4.29%  ' }
    @{ Cell = 'E23'; Value = '  -0.28%  ' }
    @{ Cell = 'E24'; Value = '  -1.94%  ' }
    @{ Cell = 'D25'; Value = '65.51' }
    @{ Cell = 'E25'; Value = '  -0.05%  ' }
    @{ Cell = 'D26'; Value = '0.993' }
    @{ Cell = 'E26'; Value = '  -0.77%  ' }
    @{ Cell = 'D27'; Value = '2.533.95' }
    @{ Cell = 'E27'; Value = '  -1.27%  ' }
    @{ Cell = 'D28'; Value = '0.155' }
    @{ Cell = 'E28'; Value = '  -1.27%  ' }
    @{ Cell = 'E29'; Value = '  -5.29%  ' }
    @{ Cell = 'D30'; Value = '7.56' }
    @{ Cell = 'E30'; Value = '  +4.07%  ' }
    @{ Cell = 'D31'; Value = '173.13' }
    @{ Cell = 'E31'; Value = '  -0.22%  ' }
    @{ Cell = 'D32'; Value = '0.0₃0731' }
    @{ Cell = 'E32'; Value = '  -0.98%  ' }
    @{ Cell = 'E33'; Value = '  -0.23%  ' }
    @{ Cell = 'D34'; Value = '6.16' }
    @{ Cell = 'E34'; Value = '  -0.26%  ' }
    @{ Cell = 'D35'; Value = '1.13' }
    @{ Cell = 'E35'; Value = '  -0.16%  ' }
    @{ Cell = 'E36'; Value = '  -0.10%  ' }
    @{ Cell = 'D37'; Value = '0.997' }
    @{ Cell = 'E37'; Value = '  -0.08%  ' }
    @{ Cell = 'D38'; Value = '17.99' }
    @{ Cell = 'E38'; Value = '  +0.88%  ' }
    @{ Cell = 'D39'; Value = '1.24' }
    @{ Cell = 'E39'; Value = '  +4.66%  ' }
    @{ Cell = 'D40'; Value = '3.84' }
    @{ Cell = 'E40'; Value = '  +1.06%  ' }
    @{ Cell = 'D41'; Value = '36.72' }
    @{ Cell = 'E41'; Value = '  +1.17%  ' }
    @{ Cell = 'B42'; Value = 'Stacks' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D42'; Value = '1.46' }
    @{ Cell = 'E42'; Value = '  +0.72%  ' }
    @{ Cell = 'B43'; Value = 'SuiNetwork' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui' }
    @{ Cell = 'D43'; Value = '0.806' }
    @{ Cell = 'E43'; Value = '  -0.64%  ' }
    @{ Cell = 'D44'; Value = '135.47' }
    @{ Cell = 'E44'; Value = '  +10.90%  ' }
    @{ Cell = 'D45'; Value = '3.39' }
    @{ Cell = 'E45'; Value = '  -0.40%  ' }
    @{ Cell = 'D46'; Value = '4.96' }
    @{ Cell = 'E46'; Value = '  +3.13%  ' }
    @{ Cell = 'D47'; Value = '256.30' }
    @{ Cell = 'E47'; Value = '  -2.10%  ' }
    @{ Cell = 'D48'; Value = '0.573' }
    @{ Cell = 'E48'; Value = '  -2.30%  ' }
    @{ Cell = 'D49'; Value = '0.0917' }
    @{ Cell = 'E49'; Value = '  -0.49%  ' }
    @{ Cell = 'D50'; Value = '0.0491' }
    @{ Cell = 'E50'; Value = '  -0.48%  ' }
    @{ Cell = 'D51'; Value = '0.0214' }
    @{ Cell = 'E51'; Value = '  +0.93%  ' }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
}
